{"js": "// Office.js (Word JavaScript API) script.\n// Applies the same edits described by the OOXML diff:\n//   1. After \"Wall height: 0.5m\" insert a blank paragraph followed by a\n//      paragraph with the text \"Skid-skid distance: = 9.47293m\".\n//   2. At the end of the document, turn the trailing empty paragraph into a\n//      page break, then append \"shakles for lifting:\" and the McMaster-Carr\n//      shackle URL as two new paragraphs.\n\n// --- 1. Insert the skid-skid distance note under \"Wall height: 0.5m\" ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet wallHeightPara = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Wall height: 0.5m\") {\n    wallHeightPara = paragraphs.items[i];\n    break;\n  }\n}\n\nif (wallHeightPara) {\n  const blankPara = wallHeightPara.insertParagraph(\"\", \"After\");\n  const distancePara = blankPara.insertParagraph(\"Skid-skid distance: = \", \"After\");\n  distancePara.insertText(\"9.47293\", \"End\");\n  await context.sync();\n  distancePara.insertText(\"m\", \"End\");\n  await context.sync();\n}\n\n// --- 2. Add the shackles-for-lifting note at the end of the document ---\nconst tailParagraphs = context.document.body.paragraphs;\ntailParagraphs.load(\"text\");\nawait context.sync();\n\nconst lastParagraph = tailParagraphs.items[tailParagraphs.items.length - 1];\nlastParagraph.getRange(\"End\").insertBreak(\"Page\", \"Start\");\nawait context.sync();\n\nconst refreshedParagraphs = context.document.body.paragraphs;\nrefreshedParagraphs.load(\"text\");\nawait context.sync();\n\nconst newLastParagraph = refreshedParagraphs.items[refreshedParagraphs.items.length - 1];\nnewLastParagraph.insertText(\"shakles for lifting:\", \"Start\");\nnewLastParagraph.insertParagraph(\n  \"https://www.mcmaster.com/products/shackles/pin-type~safety/application~for-lifting/material~stainless-steel/safety-factor~5-1/\",\n  \"After\"\n);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the same edits described by the OOXML diff:\n#   1. After \"Wall height: 0.5m\" insert a blank paragraph followed by a\n#      paragraph with the text \"Skid-skid distance: = 9.47293m\".\n#   2. At the end of the document, turn the trailing empty paragraph into a\n#      page break, then append \"shakles for lifting:\" and the McMaster-Carr\n#      shackle URL as two new paragraphs.\n\n$d = $word.ActiveDocument\n\n# --- 1. Insert the skid-skid distance note under \"Wall height: 0.5m\" ---\n$wallHeightIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -eq \"Wall height: 0.5m`r\") {\n        $wallHeightIndex = $i\n        break\n    }\n}\n\nif ($wallHeightIndex -ge 1) {\n    $wallHeightPara = $d.Paragraphs.Item($wallHeightIndex)\n    $wallHeightPara.Range.InsertParagraphAfter()\n\n    $d = $word.ActiveDocument\n    $blankPara = $d.Paragraphs.Item($wallHeightIndex + 1)\n    $blankPara.Range.InsertParagraphAfter()\n\n    $d = $word.ActiveDocument\n    $skidPara = $d.Paragraphs.Item($wallHeightIndex + 2)\n    $skidPara.Range.InsertBefore(\"Skid-skid distance: = 9.47293m\")\n}\n\n# --- 2. Add the shackles-for-lifting note at the end of the document ---\n$d = $word.ActiveDocument\n$lastIndex = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($lastIndex)\n$lastPara.Range.InsertParagraphAfter()\n\n$d = $word.ActiveDocument\n$breakPara = $d.Paragraphs.Item($lastIndex)\n$breakPara.Range.InsertBefore([char]12)\n\n$d = $word.ActiveDocument\n$shacklesIndex = $d.Paragraphs.Count\n$shacklesPara = $d.Paragraphs.Item($shacklesIndex)\n$shacklesPara.Range.InsertBefore(\"shakles for lifting:\")\n$shacklesPara.Range.InsertParagraphAfter()\n\n$d = $word.ActiveDocument\n$urlIndex = $d.Paragraphs.Count\n$urlPara = $d.Paragraphs.Item($urlIndex)\n$urlPara.Range.InsertBefore(\"https://www.mcmaster.com/products/shackles/pin-type~safety/application~for-lifting/material~stainless-steel/safety-factor~5-1/\")\n"}
